# Generate Report for Handoff
#
# Inserts two new files (2b2efc2c-... and 863139b5-...) into the
# localization-status report, ahead of the two files that were already
# present (a34f14c9-... and c033de81-...), across all three sheets:
#   Overview (summary), zh-cn (per-language detail), de-de (per-language detail)

$wb = $excel.ActiveWorkbook

$hyperColor = 15570276   # RGB(100,149,237) == FF6495ED, matches the workbook's existing HyperLink font
$dateFmt    = "yyyy-mm-dd HH:mm:ss"

# Ordered list of the four "File Name" records as they should appear after
# the edit: the two newly-handed-off files first, then the two that were
# already in the report.
$files = @(
    @{ Guid = "2b2efc2c-9a61-4264-b2b9-a9e59d77dd52"; Hash = "b541e823966cc0a0e8ac9e047ae69c78c45797aa" },
    @{ Guid = "863139b5-46f8-4819-96d8-197578fdf717"; Hash = "59d544611ec7894a8d3a4b4f8d11628b35b36deb" },
    @{ Guid = "a34f14c9-6ad5-45a7-9aa2-3728445d3e96"; Hash = "5dbd702ec3bdd42ddfbc136295b141a77a4ed2fc" },
    @{ Guid = "c033de81-2661-4ffd-95d6-1938de6ae6c8"; Hash = "7e18efceb1d0866379b36ac829637fb80432e8f6" }
)

$mdCommit  = "27bf446ec502d8f0abac8162f806ff6b5629c778"
$zhCommit  = "e6c0f9cfd7630df4e9f5510f74aca85163b8f886"
$deCommit  = "d9e92ce6556b9f1e3ef74a09c205d81b978fb32a"

function Style-Hyperlink($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $hyperColor
}

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Hyperlinks.Delete()

for ($i = 0; $i -lt $files.Count; $i++) {
    $row = $i + 2
    $guid = $files[$i].Guid

    $wsOverview.Range("A$row").Value2 = "$guid.md"
    $wsOverview.Range("B$row").Value2 = "Ready for handoff"
    $wsOverview.Range("C$row").Value2 = "Ready for handoff"
    $wsOverview.Range("D$row").Value2 = "2016-13-13 22:13:16"

    $mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/$guid.md"
    $wsOverview.Hyperlinks.Add($wsOverview.Range("A$row"), $mdUrl, $null, $null, "$guid.md") | Out-Null

    # Hyperlinks.Add re-stamps its own (theme-based) style, so (re)apply the
    # workbook's custom hyperlink look *after* adding the link.
    Style-Hyperlink $wsOverview.Range("A$row")
}

# ---------------------------------------------------------------------
# Per-language detail sheets: zh-cn / de-de
# ---------------------------------------------------------------------
$langSheets = @(
    @{ Name = "zh-cn"; Lang = "zh-cn"; Commit = $zhCommit; HandoffTime = "2016-03-13 22:13:12" },
    @{ Name = "de-de"; Lang = "de-de"; Commit = $deCommit; HandoffTime = "2016-03-13 22:13:16" }
)

foreach ($langInfo in $langSheets) {
    $ws = $wb.Worksheets.Item($langInfo.Name)
    $ws.Hyperlinks.Delete()

    for ($i = 0; $i -lt $files.Count; $i++) {
        $row = $i + 2
        $guid = $files[$i].Guid
        $hash = $files[$i].Hash
        $lang = $langInfo.Lang

        $xlfName = "$guid.$hash.$lang.xlf"

        $ws.Range("A$row").Value2 = "$guid.md"
        $ws.Range("B$row").Value2 = ".md"
        $ws.Range("C$row").Value2 = "Ready for handoff"
        $ws.Range("D$row").Value2 = $xlfName
        $ws.Range("E$row").Value2 = $langInfo.HandoffTime
        $ws.Range("H$row").Value2 = "0001-01-01 00:00:00"
        $ws.Range("I$row").Value2 = "Include"
        $ws.Range("E$row").NumberFormat = $dateFmt

        $mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/$guid.md"
        $xlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$($langInfo.Commit)/ol-handoff/OpenLocalizationTestOrg/oltest.$lang/ci/high/$xlfName"

        $ws.Hyperlinks.Add($ws.Range("A$row"), $mdUrl, $null, $null, "$guid.md") | Out-Null
        $ws.Hyperlinks.Add($ws.Range("B$row"), $mdUrl, $null, $null, ".md") | Out-Null
        $ws.Hyperlinks.Add($ws.Range("D$row"), $xlfUrl, $null, $null, $xlfName) | Out-Null

        # Hyperlinks.Add re-stamps its own (theme-based) style, so (re)apply
        # the workbook's custom hyperlink look *after* adding the links.
        Style-Hyperlink $ws.Range("A$row")
        Style-Hyperlink $ws.Range("B$row")
        Style-Hyperlink $ws.Range("D$row")
    }
}

Write-Output "Report regenerated for handoff."
